$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reports")

# --- Header block updates -------------------------------------------------
$ws.Range("B1").Value = "akhil"
$ws.Range("B2").Value = "fkjds"
# B5 holds a date-like string; format as Text first so it is not
# auto-converted into a date serial number.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2020-08-11"

# --- Existing data rows 10 & 11: update in place ---------------------------
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2020-08-08"
$ws.Range("B10").Value = "DJFLDSK"
$ws.Range("E10").Value = 0

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "2020-08-08"
$ws.Range("B11").Value = "DKJFHKDJ"
$ws.Range("E11").Value = 0

# --- Row 12 used to be the TOTAL row; it becomes a normal data row. Seed it
# with the formatting (borders/fill/font) of the existing data rows first,
# then overwrite the values.
$ws.Range("A10:E10").Copy($ws.Range("A12:E12"))
$ws.Range("A12").Value = "2020-08-10"
$ws.Range("B12").Value = "DJFLDSK"
$ws.Range("C12").Value = "IN-0003"
$ws.Range("D12").Value = "INVOICE"
$ws.Range("E12").Value = 0

# --- Row 13: brand-new data row -------------------------------------------
$ws.Range("A10:E10").Copy($ws.Range("A13:E13"))
$ws.Range("A13").Value = "2020-08-11"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "44444444444"
$ws.Range("C13").Value = "IN-0004"
$ws.Range("D13").Value = "INVOICE"
$ws.Range("E13").Value = 4012.8

# --- Row 14: new TOTAL row. Set the formula first (so its dependents,
# including the brand-new E13, stay correctly tracked), then paste the
# header row's formatting (borders/fill/font) onto it, then the label.
$ws.Range("E14").Formula = "= SUM(E10:E13)"
$ws.Range("D9:E9").Copy()
$ws.Range("D14:E14").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D14").Value = "TOTAL"
